# Daily attendance processing - 2025-10-06 06:51:38
# Applies the attendance-report updates described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" email lists (same people, new order) ---
$ws.Range("G3").Value  = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G25").Value = "Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

$ws.Range("G12").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G34").Value = "Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"

# --- Year 3 / C1 summary metrics (K/L block) ---
$ws.Range("L6").Value  = 5       # Recorded Sessions

$ws.Range("L8").Value  = 37      # Pending Sessions

# Percentage metrics are stored as literal text (e.g. "11.4%"), not numbers,
# so plain assignment (which Excel auto-parses as a percentage number) would
# change both the stored type and the cell style. Force the value in as text
# by quote-prefixing it, then fix the resulting cell style back to the
# original (General, no quote-prefix) style by pulling the number format
# from a neighbouring cell that already carries it.
$ws.Range("L9").Value = "'11.4%"   # Coverage %
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").Value = "'29.5%"  # Average Attendance %
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- Year 3 / C2 group statistics row (row 16) ---
$ws.Range("O16").Value = 3        # Recorded
$ws.Range("Q16").Value = 18       # Pending

$ws.Range("R16").Value = "'13.6%" # Coverage %
$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)

$ws.Range("S16").Value = "'23.3%" # Avg Attendance %
$ws.Range("Q16").Copy()
$ws.Range("S16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row 41 (Year 3 / C2 / PHYSIOLOGY / Session 1) newly recorded ---
$ws.Range("G41").Value = "marina_atef@med.asu.edu.eg"
$ws.Range("H41").Value = "38/246"
$ws.Range("I41").Value = "Recorded"

# Match the row's fill/formatting to the other "Recorded" rows (e.g. row 3)
# by copying formats only, so the existing "Recorded" cell style is reused
# instead of minting a new one.
$ws.Range("A3:I3").Copy()
$ws.Range("A41:I41").PasteSpecial(-4122)
$excel.CutCopyMode = $false
